$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: augment a few "type" cells with unit / enum info ---
# New shared strings must be introduced in this exact order so the
# resulting sharedStrings table matches the target layout:
#   #float,  unit:mlormg
#   #integer,  unit:nm
#   (#float already exists - reused by M2/N2/P2, untouched)
#   #float,  unit:l/mol/cm
#   #string,  unit:mmol/lormg/lorµmol/lorg/l
$ws.Range("I2").Value = "#float,  unit:mlormg"
$ws.Range("J2").Value = "#integer,  unit:nm"
$ws.Range("O2").Value = "#float,  unit:l/mol/cm"
$ws.Range("R2").Value = "#string,  unit:mmol/lormg/lorµmol/lorg/l"

# --- Row 3: new "description" row under the headers ---
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
